# Update gh-pages output data (generated at 456a3b4)
# This script updates "want-to-go" counts (column F) across the four
# worksheets, updates a couple of event time ranges (column E), flips a
# sold-out flag to text (column G on 本地生活), and appends a new event
# row (row 45) to the 展览 sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1F = @{
    5  = 188
    6  = 15
    7  = 775
    8  = 75
    9  = 9965
    11 = 3164
    13 = 2416
    14 = 2735
    15 = 1397
    17 = 2119
    19 = 89
    20 = 377
    22 = 108
    23 = 304
    24 = 264
    25 = 186
    26 = 608
    28 = 1243
    29 = 98
    32 = 2230
    33 = 2903
    34 = 9
    37 = 371
    38 = 4
    41 = 96
    44 = 35
}
foreach ($row in $sheet1F.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1F[$row]
}

# Event on row 32 had its start time pushed earlier (10:00 -> 09:30)
$ws1.Range("E32").Value = "2024.10.01 09:30-10.02 17:00"

# New event appended as row 45 (index 44 in column A)
$ws1.Range("A44").Copy($ws1.Range("A45"))
$ws1.Range("A45").Value = 44

$ws1.Cells.Item(45, 2).NumberFormat = "@"
$ws1.Cells.Item(45, 2).Value = "2024-11-23"
$ws1.Cells.Item(45, 2).ClearFormats()

$ws1.Range("C45").Value = "北京·代号鸢only同人展"
$ws1.Range("D45").Value = "北花园路1号 超级蜂巢"
$ws1.Range("E45").Value = "2024.11.23 10:00-11.23 17:00"
$ws1.Range("F45").Value = 0
$ws1.Range("G45").Value = 68
$ws1.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=90673"
$ws1.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202408/CUGDQpjZ1723553671194.jpeg"

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$sheet2F = @{
    4  = 167
    15 = 167
}
foreach ($row in $sheet2F.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2F[$row]
}

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$sheet3F = @{
    2 = 733
    3 = 966
    5 = 1936
}
foreach ($row in $sheet3F.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3F[$row]
}

# Ticket price column flips from numeric 0 to "sold out" text
$ws3.Range("G5").Value = "已售罄"

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4F = @{
    2  = 733
    3  = 966
    8  = 188
    9  = 15
    10 = 775
    11 = 75
    12 = 9965
    13 = 167
    16 = 3164
    17 = 2417
    18 = 2735
    21 = 2119
    23 = 89
    24 = 108
    25 = 304
    26 = 264
    27 = 608
    29 = 1243
    33 = 2230
    35 = 2903
    39 = 371
    41 = 4
    48 = 35
    49 = 167
}
foreach ($row in $sheet4F.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4F[$row]
}

# Same event time-range fix as on 展览!E32
$ws4.Range("E33").Value = "2024.10.01 09:30-10.02 17:00"

Write-Host "Update complete"
